# Updated cryptos list on Mon Feb 27 23:47:03 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are plain text (price strings using "." as both thousands
# and decimal separators in this export); force text format so Excel does
# not reinterpret them as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.509.27"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.633.73"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9983"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9990"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.79"
$ws.Range("E6").Value = "  -1.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.04"
$ws.Range("E8").Value = "  -1.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3648"
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("E10").Value = "  -3.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08132"
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9990"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("E13").Value = "  -2.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.598"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001252"
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.279"
$ws.Range("E16").Value = "  -2.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.626.56"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.02"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06938"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("E20").Value = "  -2.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.450"
$ws.Range("E21").Value = "  -2.23%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "23.515.91"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.79"
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.278"
$ws.Range("E25").Value = "  +5.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.434"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.33"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.59"
$ws.Range("E28").Value = "  -1.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.303"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.99"
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.311"
$ws.Range("E31").Value = "  -4.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.813.90"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.876"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.06"
$ws.Range("E34").Value = "  +5.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9649"
$ws.Range("E35").Value = "  -1.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02810"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2549"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.07234"
$ws.Range("E38").Value = "  -3.36%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.136"
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08839"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7131"
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.358"
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.31"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.40"
$ws.Range("E44").Value = "  -2.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6548"
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.347"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.010"
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08022"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "125.82"
$ws.Range("E51").Value = "  -4.09%  "
